$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text
# (matching the original inline-string cell type). We force text format,
# assign the value, then restore the default "Normal" style so no visible
# formatting change is introduced.
$numericLookingCells = @{
    'D5' = '240.25'
    'D6' = '0.6276'
    'D8' = '0.07635'
    'D10' = '24.70'
    'D12' = '5.033'
    'D13' = '0.6783'
    'D15' = '83.19'
    'D16' = '6.157'
    'D18' = '226.53'
    'D20' = '0.9999'
    'D21' = '7.498'
    'D22' = '0.9995'
    'D23' = '158.16'
    'D24' = '0.1381'
    'D25' = '8.403'
    'D27' = '1.386'
    'D28' = '1.459'
    'D29' = '0.05609'
    'D31' = '4.076'
    'D34' = '0.6909'
    'D36' = '0.01804'
    'D38' = '2.720'
    'D39' = '6.381'
    'D40' = '0.9048'
    'D42' = '101.56'
    'D43' = '66.03'
    'D44' = '7.188'
    'D45' = '0.00000000119'
    'D46' = '0.4011'
    'D47' = '9.024'
    'D48' = '1.677'
    'D49' = '0.1140'
    'D50' = '0.05705'
    'D51' = '0.4627'
}

foreach ($addr in $numericLookingCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $numericLookingCells.Keys) {
    $ws.Range($addr).Value = $numericLookingCells[$addr]
}
foreach ($addr in $numericLookingCells.Keys) {
    $ws.Range($addr).Style = "Normal"
}

# Plain text cells (coin names, links, non-numeric-looking prices, and
# all percentage strings) can be assigned directly.
$textCells = @{
    'D2' = '29.386.46'
    'E2' = '  -0.03%  '
    'D3' = '1.848.02'
    'E3' = '  -0.02%  '
    'E5' = '  -0.03%  '
    'E6' = '  -0.36%  '
    'E7' = '  -0.01%  '
    'E8' = '  +0.27%  '
    'E9' = '  -1.03%  '
    'E10' = '  +0.89%  '
    'E12' = '  +0.54%  '
    'E13' = '  -0.13%  '
    'E14' = '  -1.32%  '
    'E15' = '  -0.60%  '
    'E16' = '  -0.03%  '
    'D17' = '29.404.86'
    'E18' = '  -1.05%  '
    'E21' = '  +0.67%  '
    'E22' = '  -0.10%  '
    'E23' = '  +0.50%  '
    'E24' = '  -0.61%  '
    'E25' = '  +0.17%  '
    'E26' = '  +0.22%  '
    'E27' = '  +5.42%  '
    'E28' = '  -0.49%  '
    'E29' = '  -0.12%  '
    'E30' = '  +0.23%  '
    'E31' = '  +1.04%  '
    'E32' = '  -0.70%  '
    'E33' = '  +0.49%  '
    'E34' = '  -2.61%  '
    'E35' = '  -0.20%  '
    'B36' = 'VeChain'
    'C36' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E36' = '  +0.14%  '
    'B37' = 'Maker'
    'C37' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D37' = '1.230.30'
    'E37' = '  -0.38%  '
    'E38' = '  -2.01%  '
    'E39' = '  -1.36%  '
    'E40' = '  -0.39%  '
    'E41' = '  +0.04%  '
    'E42' = '  +0.09%  '
    'E43' = '  -0.09%  '
    'E44' = '  +0.19%  '
    'E45' = '  -1.92%  '
    'E46' = '  -0.14%  '
    'E47' = '  +0.22%  '
    'E48' = '  -0.33%  '
    'E49' = '  +1.54%  '
    'E50' = '  -0.10%  '
    'E51' = '  +0.05%  '
}

foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}
